$d = $word.ActiveDocument
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($p.Range.Text -notmatch "stop the exception propagation") {
    throw "Unexpected last paragraph text: $($p.Range.Text)"
}
$rng = $p.Range
$xml = @'
<w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>If you want to stop the exception propagation then use @Around advice.</w:t></w:r></w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>43</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> AOP After </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">(finally) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Advice Type</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">@After runs after a method is completed. Regardless of outcome / exceptions (works like “finally” block). After will execute before </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>AfterThrowing</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Use cases:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Log the exception and/or perform auditing</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Code to run regardless of method outcome</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Encapsulate this functionality in AOP aspect for easy reuse</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>After advice does not have access to the exception. If you need exception, then use @</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>AfterThrowing</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> advice</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>

'@
$rng.InsertXML($xml)
